$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 22737264
$ws.Range("I33").Value = 50001980
$ws.Range("K33").Value = 50001980
$ws.Range("M33").Value = -50001751
$ws.Range("H40").Value = 4767.467
$ws.Range("I40").Value = 4000
$ws.Range("J40").Value = 4822.2856
$ws.Range("K40").Value = 4000
$ws.Range("L40").Value = 4822.2856
$ws.Range("M40").Value = -3825
$ws.Range("N40").Value = -5172.2856
$ws.Range("H62").Value = 16740.8
$ws.Range("I62").Value = 4852.5
$ws.Range("J62").Value = 24666.334
$ws.Range("K62").Value = 4852.5
$ws.Range("L62").Value = 24666.334
$ws.Range("M62").Value = -4228.5
$ws.Range("N62").Value = -25914.334
$ws.Range("H64").Value = 7000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 7000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 7000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -7496
$ws.Range("H65").Value = 16740.8
$ws.Range("I65").Value = 4852.5
$ws.Range("J65").Value = 24666.334
$ws.Range("K65").Value = 24262.5
$ws.Range("L65").Value = 123331.67
$ws.Range("M65").Value = -21142.5
$ws.Range("N65").Value = -129571.67
$ws.Range("H67").Value = 7000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 7000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 7000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -8716
$ws.Range("H104").Value = 90
$ws.Range("I104").Value = 90
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 270
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = 1477
$ws.Range("N104").ClearContents()
$ws.Range("H106").Value = 6955.7144
$ws.Range("I106").Value = 3258.2
$ws.Range("J106").Value = 16199.5
$ws.Range("K106").Value = 3258.2
$ws.Range("L106").Value = 16199.5
$ws.Range("M106").Value = -2627.2
$ws.Range("N106").Value = -17461.5
$ws.Range("H112").Value = 4578
$ws.Range("J112").Value = 4920
$ws.Range("L112").Value = 14760
$ws.Range("N112").Value = -16976
$ws.Range("H113").Value = 15748.571
$ws.Range("I113").Value = 20387
$ws.Range("K113").Value = 20387
$ws.Range("M113").Value = -17133
$ws.Range("H121").Value = 4924
$ws.Range("J121").Value = 4924
$ws.Range("L121").Value = 14772
$ws.Range("N121").Value = -18266
$ws.Range("H132").Value = 14939.174
$ws.Range("I132").Value = 11361.667
$ws.Range("K132").Value = 34085.001
$ws.Range("M132").Value = -31555.001
$ws.Range("H137").Value = 9632.5
$ws.Range("I137").Value = 2737.261
$ws.Range("K137").Value = 8211.783
$ws.Range("M137").Value = -5661.782999999999
$ws.Range("H138").Value = 3364.1177
$ws.Range("I138").Value = 2955.077
$ws.Range("J138").Value = 4693.5
$ws.Range("K138").Value = 8865.231
$ws.Range("L138").Value = 14080.5
$ws.Range("M138").Value = -3725.231
$ws.Range("N138").Value = -24360.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 4000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -4288
$ws.Range("H11").Value = 14000
$ws.Range("J11").Value = 14000
$ws.Range("L11").Value = 14000
$ws.Range("N11").Value = -14288
$ws.Range("H17").Value = 892.8571
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 958.3333
$ws.Range("K17").Value = 500
$ws.Range("L17").Value = 958.3333
$ws.Range("M17").Value = -327
$ws.Range("N17").Value = -1304.3333
$ws.Range("H32").Value = 19227.777
$ws.Range("I32").Value = 12442.143
$ws.Range("K32").Value = 12442.143
$ws.Range("M32").Value = -12155.143
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H61").Value = 17713.291
$ws.Range("I61").Value = 3370.6365
$ws.Range("K61").Value = 3370.6365
$ws.Range("M61").Value = -3158.6365
$ws.Range("H63").Value = 2713.375
$ws.Range("I63").Value = 1979.8
$ws.Range("K63").Value = 1979.8
$ws.Range("M63").Value = -1293.8
$ws.Range("H66").Value = 2713.375
$ws.Range("I66").Value = 1979.8
$ws.Range("K66").Value = 9899
$ws.Range("M66").Value = -6467
$ws.Range("H74").Value = 20095.084
$ws.Range("I74").Value = 3802.6
$ws.Range("J74").Value = 31732.572
$ws.Range("K74").Value = 3802.6
$ws.Range("L74").Value = 31732.572
$ws.Range("M74").Value = -2928.6
$ws.Range("N74").Value = -33480.572
$ws.Range("H77").Value = 20095.084
$ws.Range("I77").Value = 3802.6
$ws.Range("J77").Value = 31732.572
$ws.Range("K77").Value = 19013
$ws.Range("L77").Value = 158662.86
$ws.Range("M77").Value = -14645
$ws.Range("N77").Value = -167398.86
$ws.Range("H122").Value = 2965576.8
$ws.Range("I122").Value = 4608675
$ws.Range("K122").Value = 13826025
$ws.Range("M122").Value = -13823575
$ws.Range("H132").Value = 2954541.2
$ws.Range("I132").Value = 4683.2173
$ws.Range("K132").Value = 14049.6519
$ws.Range("M132").Value = -11519.6519
$ws.Range("H136").Value = 17713.291
$ws.Range("I136").Value = 3370.6365
$ws.Range("K136").Value = 10111.9095
$ws.Range("M136").Value = -7561.9095

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 784.9
$ws.Range("I12").Value = 310.4
$ws.Range("J12").Value = 1259.4
$ws.Range("K12").Value = 310.4
$ws.Range("L12").Value = 1259.4
$ws.Range("M12").Value = -142.4
$ws.Range("N12").Value = -1595.4
$ws.Range("H17").Value = 750
$ws.Range("J17").Value = 750
$ws.Range("L17").Value = 750
$ws.Range("N17").Value = -1094
$ws.Range("H20").Value = 15192170
$ws.Range("I20").Value = 47634390
$ws.Range("K20").Value = 47634390
$ws.Range("M20").Value = -47634143
$ws.Range("H22").Value = 35714508
$ws.Range("I22").Value = 47619180
$ws.Range("K22").Value = 47619180
$ws.Range("M22").Value = -47619007
$ws.Range("H87").Value = 250000
$ws.Range("J87").Value = 250000
$ws.Range("L87").Value = 250000
$ws.Range("N87").Value = -252496
$ws.Range("H90").Value = 250000
$ws.Range("J90").Value = 250000
$ws.Range("L90").Value = 750000
$ws.Range("N90").Value = -762480
$ws.Range("H107").Value = 2275.4614
$ws.Range("I107").Value = 2200
$ws.Range("J107").Value = 2340.1428
$ws.Range("K107").Value = 2200
$ws.Range("L107").Value = 2340.1428
$ws.Range("M107").Value = -280
$ws.Range("N107").Value = -6180.1428
$ws.Range("H132").Value = 71198.14
$ws.Range("I132").Value = 40000
$ws.Range("J132").Value = 76397.836
$ws.Range("K132").Value = 40000
$ws.Range("L132").Value = 76397.836
$ws.Range("M132").Value = -34940
$ws.Range("N132").Value = -86517.836
$ws.Range("H134").Value = 14731.315
$ws.Range("I134").Value = 4698.0586
$ws.Range("K134").Value = 14094.1758
$ws.Range("M134").Value = -11559.1758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 178.66667
$ws.Range("I5").Value = 194.4
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 194.4
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -82.4
$ws.Range("N5").Value = -324
$ws.Range("H12").Value = 1077.6
$ws.Range("I12").Value = 378.8
$ws.Range("J12").Value = 3174
$ws.Range("K12").Value = 378.8
$ws.Range("L12").Value = 3174
$ws.Range("M12").Value = -208.8
$ws.Range("N12").Value = -3514
$ws.Range("H16").Value = 5851.091
$ws.Range("I16").Value = 1058.3334
$ws.Range("J16").Value = 11602.4
$ws.Range("K16").Value = 1058.3334
$ws.Range("L16").Value = 11602.4
$ws.Range("M16").Value = -771.3334
$ws.Range("N16").Value = -12176.4
$ws.Range("H25").Value = 2275.077
$ws.Range("I25").Value = 1816.6666
$ws.Range("K25").Value = 1816.6666
$ws.Range("M25").Value = -1642.6666
$ws.Range("H31").Value = 16221.105
$ws.Range("I31").Value = 8920.4375
$ws.Range("J31").Value = 21530.682
$ws.Range("K31").Value = 8920.4375
$ws.Range("L31").Value = 21530.682
$ws.Range("M31").Value = -8625.4375
$ws.Range("N31").Value = -22120.682
$ws.Range("H34").Value = 16221.105
$ws.Range("I34").Value = 8920.4375
$ws.Range("J34").Value = 21530.682
$ws.Range("K34").Value = 8920.4375
$ws.Range("L34").Value = 21530.682
$ws.Range("M34").Value = -8718.4375
$ws.Range("N34").Value = -21934.682
$ws.Range("H58").Value = 13768.667
$ws.Range("I58").Value = 6514
$ws.Range("J58").Value = 20259.684
$ws.Range("K58").Value = 6514
$ws.Range("L58").Value = 20259.684
$ws.Range("M58").Value = -6311
$ws.Range("N58").Value = -20665.684
$ws.Range("H60").Value = 13378.667
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H99").Value = 3446178
$ws.Range("I99").Value = 5937784
$ws.Range("K99").Value = 5937784
$ws.Range("M99").Value = -5936286
$ws.Range("H107").Value = 2691.2903
$ws.Range("I107").Value = 1250.5264
$ws.Range("J107").Value = 4972.5
$ws.Range("K107").Value = 1250.5264
$ws.Range("L107").Value = 4972.5
$ws.Range("M107").Value = 669.4736
$ws.Range("N107").Value = -8812.5
$ws.Range("H113").Value = 5851.091
$ws.Range("I113").Value = 1058.3334
$ws.Range("J113").Value = 11602.4
$ws.Range("K113").Value = 1058.3334
$ws.Range("L113").Value = 11602.4
$ws.Range("M113").Value = 1111.6666
$ws.Range("N113").Value = -15942.4
$ws.Range("H126").Value = 3446178
$ws.Range("I126").Value = 5937784
$ws.Range("K126").Value = 17813352
$ws.Range("M126").Value = -17810882
$ws.Range("H132").Value = 16675.6
$ws.Range("I132").Value = 8376
$ws.Range("K132").Value = 25128
$ws.Range("M132").Value = -22598
$ws.Range("H136").Value = 13768.667
$ws.Range("I136").Value = 6514
$ws.Range("J136").Value = 20259.684
$ws.Range("K136").Value = 19542
$ws.Range("L136").Value = 60779.052
$ws.Range("M136").Value = -16992
$ws.Range("N136").Value = -65879.052

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 194.66667
$ws.Range("I2").Value = 232.3
$ws.Range("K2").Value = 1393.8
$ws.Range("M2").Value = -1280.8
$ws.Range("H5").Value = 3050889
$ws.Range("I5").Value = 2299.6667
$ws.Range("J5").Value = 4880042.5
$ws.Range("K5").Value = 6899.000100000001
$ws.Range("L5").Value = 14640127.5
$ws.Range("M5").Value = -6787.000100000001
$ws.Range("N5").Value = -14640351.5
$ws.Range("H12").Value = 200.8
$ws.Range("J12").Value = 238.5
$ws.Range("L12").Value = 715.5
$ws.Range("N12").Value = -1061.5
$ws.Range("H22").Value = 772.3333
$ws.Range("I22").Value = 340.6154
$ws.Range("J22").Value = 1894.8
$ws.Range("K22").Value = 1021.8462
$ws.Range("L22").Value = 5684.4
$ws.Range("M22").Value = -852.8462000000001
$ws.Range("N22").Value = -6022.4
$ws.Range("H27").Value = 772.3333
$ws.Range("I27").Value = 340.6154
$ws.Range("J27").Value = 1894.8
$ws.Range("K27").Value = 1021.8462
$ws.Range("L27").Value = 5684.4
$ws.Range("M27").Value = -919.8462000000001
$ws.Range("N27").Value = -5888.4
$ws.Range("H38").Value = 52.333332
$ws.Range("I38").Value = 9.5
$ws.Range("K38").Value = 28.5
$ws.Range("M38").Value = 318.5
$ws.Range("H39").Value = 11999.5
$ws.Range("J39").Value = 11999.5
$ws.Range("L39").Value = 35998.5
$ws.Range("N39").Value = -36586.5
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312
$ws.Range("H51").Value = 2600
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 15000
$ws.Range("N51").Value = -15920
$ws.Range("H80").Value = 17807.54
$ws.Range("I80").Value = 7116.6665
$ws.Range("J80").Value = 26971.143
$ws.Range("K80").Value = 21349.9995
$ws.Range("L80").Value = 80913.429
$ws.Range("M80").Value = -20413.9995
$ws.Range("N80").Value = -82785.429
$ws.Range("H83").Value = 17807.54
$ws.Range("I83").Value = 7116.6665
$ws.Range("J83").Value = 26971.143
$ws.Range("K83").Value = 64049.9985
$ws.Range("L83").Value = 242740.287
$ws.Range("M83").Value = -59369.9985
$ws.Range("N83").Value = -252100.287
$ws.Range("H98").Value = 4441.778
$ws.Range("I98").Value = 1892
$ws.Range("J98").Value = 4760.5
$ws.Range("K98").Value = 5676
$ws.Range("L98").Value = 14281.5
$ws.Range("M98").Value = -4178
$ws.Range("N98").Value = -17277.5
$ws.Range("H104").Value = 828778.44
$ws.Range("I104").Value = 6374.75
$ws.Range("J104").Value = 1768668.4
$ws.Range("K104").Value = 19124.25
$ws.Range("L104").Value = 5306005.199999999
$ws.Range("M104").Value = -16503.25
$ws.Range("N104").Value = -5311247.199999999
$ws.Range("H106").Value = 4333.1665
$ws.Range("J106").Value = 5749.75
$ws.Range("L106").Value = 17249.25
$ws.Range("N106").Value = -19141.25
$ws.Range("H114").Value = 3850
$ws.Range("I114").Value = 3850
$ws.Range("K114").Value = 11550
$ws.Range("M114").Value = -8296
$ws.Range("H131").Value = 1469.73
$ws.Range("I131").Value = 1075.2858
$ws.Range("J131").Value = 1499.4193
$ws.Range("K131").Value = 3225.8574
$ws.Range("L131").Value = 4498.257900000001
$ws.Range("M131").Value = 1814.1426
$ws.Range("N131").Value = -14578.2579
$ws.Range("H135").Value = 3050889
$ws.Range("I135").Value = 2299.6667
$ws.Range("J135").Value = 4880042.5
$ws.Range("K135").Value = 20697.0003
$ws.Range("L135").Value = 43920382.5
$ws.Range("M135").Value = -18162.0003
$ws.Range("N135").Value = -43925452.5
$ws.Range("H139").Value = 16156.571
$ws.Range("I139").Value = 16156.571
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 48469.713
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -43329.713
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 21538.5
$ws.Range("H65").Value = 21538.5
$ws.Range("H80").Value = 20323
$ws.Range("I80").Value = 15519.1875
$ws.Range("K80").Value = 15519.1875
$ws.Range("M80").Value = -14521.1875
$ws.Range("H83").Value = 20323
$ws.Range("I83").Value = 15519.1875
$ws.Range("K83").Value = 77595.9375
$ws.Range("M83").Value = -72603.9375
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H102").Value = 10403244
$ws.Range("I102").Value = 12289652
$ws.Range("K102").Value = 12289652
$ws.Range("M102").Value = -12288030
$ws.Range("H113").Value = 3770.8
$ws.Range("I113").Value = 3213.5
$ws.Range("K113").Value = 3213.5
$ws.Range("M113").Value = -1043.5
$ws.Range("H122").Value = 1317300.9
$ws.Range("I122").Value = 2131971
$ws.Range("J122").Value = 13828.7
$ws.Range("K122").Value = 6395913
$ws.Range("L122").Value = 41486.10000000001
$ws.Range("M122").Value = -6393463
$ws.Range("N122").Value = -46386.10000000001
$ws.Range("H132").Value = 21498.666
$ws.Range("I132").Value = 12871.25
$ws.Range("J132").Value = 38753.5
$ws.Range("K132").Value = 38613.75
$ws.Range("L132").Value = 116260.5
$ws.Range("M132").Value = -36083.75
$ws.Range("N132").Value = -121320.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 32262002
$ws.Range("I22").Value = 2243.6086
$ws.Range("J22").Value = 125008810
$ws.Range("K22").Value = 2243.6086
$ws.Range("L22").Value = 125008810
$ws.Range("M22").Value = -1948.6086
$ws.Range("N22").Value = -125009400
$ws.Range("H27").Value = 32262002
$ws.Range("I27").Value = 2243.6086
$ws.Range("J27").Value = 125008810
$ws.Range("K27").Value = 2243.6086
$ws.Range("L27").Value = 125008810
$ws.Range("M27").Value = -2136.6086
$ws.Range("N27").Value = -125009024
$ws.Range("H40").Value = 2564792.5
$ws.Range("I40").Value = 3202.111
$ws.Range("J40").Value = 4211529
$ws.Range("K40").Value = 3202.111
$ws.Range("L40").Value = 4211529
$ws.Range("M40").Value = -3066.111
$ws.Range("N40").Value = -4211801
$ws.Range("H46").Value = 2754.9443
$ws.Range("I46").Value = 1566.3334
$ws.Range("J46").Value = 3943.5557
$ws.Range("K46").Value = 1566.3334
$ws.Range("L46").Value = 3943.5557
$ws.Range("M46").Value = -1378.3334
$ws.Range("N46").Value = -4319.5557
$ws.Range("H63").Value = 20500
$ws.Range("J63").Value = 20500
$ws.Range("L63").Value = 20500
$ws.Range("N63").Value = -21998
$ws.Range("H66").Value = 20500
$ws.Range("J66").Value = 20500
$ws.Range("L66").Value = 61500
$ws.Range("N66").Value = -68988
$ws.Range("H68").Value = 5868
$ws.Range("I68").Value = 3601.4
$ws.Range("K68").Value = 3601.4
$ws.Range("M68").Value = -2852.4
$ws.Range("H71").Value = 5868
$ws.Range("I71").Value = 3601.4
$ws.Range("K71").Value = 18007
$ws.Range("M71").Value = -14263
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 26627140
$ws.Range("I122").Value = 41662960
$ws.Range("J122").Value = 3133667.2
$ws.Range("K122").Value = 124988880
$ws.Range("L122").Value = 9401001.600000001
$ws.Range("M122").Value = -124986430
$ws.Range("N122").Value = -9405901.600000001
$ws.Range("H132").Value = 4028701
$ws.Range("I132").Value = 40000
$ws.Range("J132").Value = 4471890
$ws.Range("K132").Value = 120000
$ws.Range("L132").Value = 13415670
$ws.Range("M132").Value = -117470
$ws.Range("N132").Value = -13420730
$ws.Range("H136").Value = 15807.151
$ws.Range("I136").Value = 23186
$ws.Range("K136").Value = 69558
$ws.Range("M136").Value = -67008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 7925
$ws.Range("I21").Value = 1769.2307
$ws.Range("J21").Value = 19357.143
$ws.Range("K21").Value = 1769.2307
$ws.Range("L21").Value = 19357.143
$ws.Range("M21").Value = -1534.2307
$ws.Range("N21").Value = -19827.143
$ws.Range("H35").Value = 7925
$ws.Range("I35").Value = 1769.2307
$ws.Range("J35").Value = 19357.143
$ws.Range("K35").Value = 1769.2307
$ws.Range("L35").Value = 19357.143
$ws.Range("M35").Value = -1479.2307
$ws.Range("N35").Value = -19937.143
$ws.Range("H45").Value = 70006.5
$ws.Range("J45").Value = 70006.5
$ws.Range("L45").Value = 70006.5
$ws.Range("N45").Value = -70988.5
$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
$ws.Range("H81").Value = 7000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 7000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H122").Value = 430554.06
$ws.Range("I122").Value = 586104.6
$ws.Range("J122").Value = 8345.429
$ws.Range("K122").Value = 1758313.8
$ws.Range("L122").Value = 25036.287
$ws.Range("M122").Value = -1755863.8
$ws.Range("N122").Value = -29936.287
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 9815.414
$ws.Range("I132").Value = 3508.5293
$ws.Range("J132").Value = 18750.166
$ws.Range("K132").Value = 10525.5879
$ws.Range("L132").Value = 56250.49800000001
$ws.Range("M132").Value = -7995.5879
$ws.Range("N132").Value = -61310.49800000001
$ws.Range("H136").Value = 9649.679
$ws.Range("I136").Value = 1409.1
$ws.Range("K136").Value = 4227.299999999999
$ws.Range("M136").Value = -1677.299999999999
